$wb = $excel.ActiveWorkbook

# Overview sheet: G2 timestamp update
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 11:16:30"

# zh-cn sheet: H2 and K2 timestamp updates
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 11:16:16"
$wsZhCn.Range("K2").Value = "2016-08-30 11:17:11"

# de-de sheet: H2 and K2 timestamp updates
# Note: de-de!H2 shared the same underlying string as Overview!G2
# ("2016-08-30 11:14:48"), so it must be updated to the same new value.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 11:16:30"
$wsDeDe.Range("K2").Value = "2016-08-30 11:17:19"
